# Individual sheet: remove the single "Trottenberg et al. (1)" data row (old row 109),
# and rename the remaining "Trottenberg et al. (2)" rows to the merged label
# "Trottenberg et al." (rows shift up by one, becoming rows 109-113).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual")
$ws.Rows(109).Delete()

for ($r = 109; $r -le 113; $r++) {
    $ws.Cells.Item($r, 1).Value = "Trottenberg et al."
}

# means sheet: remove the corresponding summary row for "Trottenberg et al. (1)"
# (old row 14), and rename the remaining summary row to "Trottenberg et al."
# (it shifts up to become row 14).
$ws2 = $wb.Worksheets.Item("means")
$ws2.Rows(14).Delete()
$ws2.Cells.Item(14, 1).Value = "Trottenberg et al."

# Restore the selections on each sheet to match the saved view state.
$ws2.Select()
$ws2.Range("A15").Select()

$ws.Select()
$ws.Range("J110").Select()
